$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows (bottom-up, using original row numbers) to expand each
# 2-row (adapt/base) group into a 3-row (adapt/adaptV1/base) group
$ws.Rows(19).Insert()
$ws.Rows(17).Insert()
$ws.Rows(15).Insert()
$ws.Rows(13).Insert()
$ws.Rows(11).Insert()
$ws.Rows(9).Insert()
$ws.Rows(7).Insert()
$ws.Rows(5).Insert()

# Copy cell formatting (style) from row 4 into each newly inserted row
$ws.Range("A4:D4").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Range("A17:D17").PasteSpecial(-4122)
$ws.Range("A20:D20").PasteSpecial(-4122)
$ws.Range("A23:D23").PasteSpecial(-4122)
$ws.Range("A26:D26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the updated table values (rows 4-27)
# Row 4
$ws.Range("A4").Value = "dataset_A"
$ws.Range("B4").Value = "chan"
$ws.Range("C4").Value = "no"
$ws.Range("D4").Value = "adapt"
$ws.Range("E4").Value = 0.506
$ws.Range("F4").Value = 0.431
$ws.Range("G4").Value = 0.456
$ws.Range("H4").Value = 0.469
$ws.Range("I4").Value = 0.419

# Row 5
$ws.Range("D5").Value = "adaptV1"
$ws.Range("E5").Value = 0.587
$ws.Range("F5").Value = 0.569
$ws.Range("G5").Value = 0.55
$ws.Range("H5").Value = 0.487
$ws.Range("I5").Value = 0.469

# Row 6
$ws.Range("D6").Value = "base"
$ws.Range("E6").Value = 0.275
$ws.Range("F6").Value = 0.344
$ws.Range("G6").Value = 0.356
$ws.Range("H6").Value = 0.269
$ws.Range("I6").Value = 0.312

# Row 7
$ws.Range("C7").Value = "temp"
$ws.Range("D7").Value = "adapt"
$ws.Range("E7").Value = 0.375
$ws.Range("F7").Value = 0.513
$ws.Range("G7").Value = 0.331
$ws.Range("H7").Value = 0.444
$ws.Range("I7").Value = 0.462

# Row 8
$ws.Range("D8").Value = "adaptV1"
$ws.Range("E8").Value = 0.594
$ws.Range("F8").Value = 0.587
$ws.Range("G8").Value = 0.531
$ws.Range("H8").Value = 0.55
$ws.Range("I8").Value = 0.419

# Row 9
$ws.Range("D9").Value = "base"
$ws.Range("E9").Value = 0.281
$ws.Range("F9").Value = 0.381
$ws.Range("G9").Value = 0.337
$ws.Range("H9").Value = 0.294
$ws.Range("I9").Value = 0.25

# Row 10
$ws.Range("B10").Value = "no"
$ws.Range("C10").Value = "no"
$ws.Range("D10").Value = "adapt"
$ws.Range("E10").Value = 0.475
$ws.Range("F10").Value = 0.475
$ws.Range("G10").Value = 0.344
$ws.Range("H10").Value = 0.506
$ws.Range("I10").Value = 0.438

# Row 11
$ws.Range("D11").Value = "adaptV1"
$ws.Range("E11").Value = 0.225
$ws.Range("F11").Value = 0.375
$ws.Range("G11").Value = 0.356
$ws.Range("H11").Value = 0.25
$ws.Range("I11").Value = 0.325

# Row 12
$ws.Range("D12").Value = "base"
$ws.Range("E12").Value = 0.231
$ws.Range("F12").Value = 0.3
$ws.Range("G12").Value = 0.306
$ws.Range("H12").Value = 0.225
$ws.Range("I12").Value = 0.2

# Row 13
$ws.Range("C13").Value = "temp"
$ws.Range("D13").Value = "adapt"
$ws.Range("E13").Value = 0.481
$ws.Range("F13").Value = 0.487
$ws.Range("G13").Value = 0.431
$ws.Range("H13").Value = 0.494
$ws.Range("I13").Value = 0.475

# Row 14
$ws.Range("D14").Value = "adaptV1"
$ws.Range("E14").Value = 0.3
$ws.Range("F14").Value = 0.469
$ws.Range("G14").Value = 0.356
$ws.Range("H14").Value = 0.25
$ws.Range("I14").Value = 0.337

# Row 15
$ws.Range("D15").Value = "base"
$ws.Range("E15").Value = 0.337
$ws.Range("F15").Value = 0.419
$ws.Range("G15").Value = 0.381
$ws.Range("H15").Value = 0.569
$ws.Range("I15").Value = 0.319

# Row 16
$ws.Range("A16").Value = "dataset_B"
$ws.Range("B16").Value = "chan"
$ws.Range("C16").Value = "no"
$ws.Range("D16").Value = "adapt"
$ws.Range("E16").Value = 0.497
$ws.Range("F16").Value = 0.521
$ws.Range("G16").Value = 0.497
$ws.Range("H16").Value = 0.538
$ws.Range("I16").Value = 0.573

# Row 17
$ws.Range("D17").Value = "adaptV1"
$ws.Range("E17").Value = 0.535
$ws.Range("F17").Value = 0.542
$ws.Range("G17").Value = 0.604
$ws.Range("H17").Value = 0.618
$ws.Range("I17").Value = 0.583

# Row 18
$ws.Range("D18").Value = "base"
$ws.Range("E18").Value = 0.37
$ws.Range("F18").Value = 0.344
$ws.Range("G18").Value = 0.396
$ws.Range("H18").Value = 0.255
$ws.Range("I18").Value = 0.406

# Row 19
$ws.Range("C19").Value = "temp"
$ws.Range("D19").Value = "adapt"
$ws.Range("E19").Value = 0.542
$ws.Range("F19").Value = 0.535
$ws.Range("G19").Value = 0.566
$ws.Range("H19").Value = 0.59
$ws.Range("I19").Value = 0.573

# Row 20
$ws.Range("D20").Value = "adaptV1"
$ws.Range("E20").Value = 0.538
$ws.Range("F20").Value = 0.608
$ws.Range("G20").Value = 0.594
$ws.Range("H20").Value = 0.615
$ws.Range("I20").Value = 0.611

# Row 21
$ws.Range("D21").Value = "base"
$ws.Range("E21").Value = 0.458
$ws.Range("F21").Value = 0.448
$ws.Range("G21").Value = 0.573
$ws.Range("H21").Value = 0.651
$ws.Range("I21").Value = 0.542

# Row 22
$ws.Range("B22").Value = "no"
$ws.Range("C22").Value = "no"
$ws.Range("D22").Value = "adapt"
$ws.Range("E22").Value = 0.49
$ws.Range("F22").Value = 0.521
$ws.Range("G22").Value = 0.601
$ws.Range("H22").Value = 0.601
$ws.Range("I22").Value = 0.559

# Row 23
$ws.Range("D23").Value = "adaptV1"
$ws.Range("E23").Value = 0.521
$ws.Range("F23").Value = 0.517
$ws.Range("G23").Value = 0.59
$ws.Range("H23").Value = 0.472
$ws.Range("I23").Value = 0.59

# Row 24
$ws.Range("D24").Value = "base"
$ws.Range("E24").Value = 0.401
$ws.Range("F24").Value = 0.318
$ws.Range("G24").Value = 0.37
$ws.Range("H24").Value = 0.385
$ws.Range("I24").Value = 0.458

# Row 25
$ws.Range("C25").Value = "temp"
$ws.Range("D25").Value = "adapt"
$ws.Range("E25").Value = 0.479
$ws.Range("F25").Value = 0.514
$ws.Range("G25").Value = 0.625
$ws.Range("H25").Value = 0.587
$ws.Range("I25").Value = 0.587

# Row 26
$ws.Range("D26").Value = "adaptV1"
$ws.Range("E26").Value = 0.521
$ws.Range("F26").Value = 0.486
$ws.Range("G26").Value = 0.51
$ws.Range("H26").Value = 0.399
$ws.Range("I26").Value = 0.556

# Row 27
$ws.Range("D27").Value = "base"
$ws.Range("E27").Value = 0.526
$ws.Range("F27").Value = 0.437
$ws.Range("G27").Value = 0.443
$ws.Range("H27").Value = 0.542
$ws.Range("I27").Value = 0.589
